$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear Factory/Dealer columns (E2:L2), populate Retailer columns (M2:P2)
$ws.Range("E2:L2").ClearContents()
$ws.Range("M2").Value = "Retailer_Reached"
$ws.Range("N2").Value = "'3"
$ws.Range("O2").Value = "2025-04-07 17:30:18"
$ws.Range("P2").Value = "Location not available"

# Row 3: clear Factory columns (E3:H3)
$ws.Range("E3:H3").ClearContents()
